$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Triggers")
$ws.Name = "Processors"
$ws.Activate() | Out-Null
$ws.Range("E58").Select() | Out-Null
